# Updated remaining queries for C3DC
# Fix the LEFT JOIN conditions used throughout the saved SQL queries on Sheet1
# (cells C2, B2, B3, B4, B5, B6, B7) so they join on the proper *_id columns
# instead of the generic "id" column:
#   std.id               -> std.study_id
#   prt."study.id"       -> prt."study.study_id"
#   prt.id                -> prt.participant_id
#   dgn."participant.id" -> dgn."participant.participant_id"   (and trt/trr/srv)
#   std.id = rfs."study.id" -> std.study_id = rfs."study.study_id"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$queryRange = $ws.Range("B2:C7")

$queryRange.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
$queryRange.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
$queryRange.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
$queryRange.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
$queryRange.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
$queryRange.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

# Move the active selection from C7 to B2
$ws.Range("B2").Select()

# Widen column C to fit the updated query text (drop the old best-fit autosize)
$ws.Columns("C").ColumnWidth = 69.83333333333333
